$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 506
$ws.Range("F7").Value = 1159
$ws.Range("F10").Value = 723
$ws.Range("F15").Value = 1618
$ws.Range("F16").Value = 17
$ws.Range("F21").Value = 1082
$ws.Range("F22").Value = 1513
$ws.Range("F23").Value = 757
$ws.Range("F24").Value = 622
$ws.Range("F25").Value = 498
$ws.Range("F29").Value = 1150
$ws.Range("F32").Value = 280
$ws.Range("F33").Value = 1368
$ws.Range("F36").Value = 3963

# Sheet: 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1038
$ws.Range("F17").Value = 26
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 258
$ws.Range("F23").Value = 121

# Sheet: 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1270
$ws.Range("F7").Value = 1010

# Sheet: 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1270
$ws.Range("F6").Value = 1010
$ws.Range("F7").Value = 507
$ws.Range("F12").Value = 1159
$ws.Range("F17").Value = 723
$ws.Range("F25").Value = 1618
$ws.Range("F26").Value = 17
$ws.Range("F30").Value = 1082
$ws.Range("F31").Value = 1513
$ws.Range("F32").Value = 757
$ws.Range("F33").Value = 622
$ws.Range("F34").Value = 498
$ws.Range("F38").Value = 258
$ws.Range("F41").Value = 1150
$ws.Range("F48").Value = 1368
$ws.Range("F50").Value = 3963
